# Auto-generated script applying scheduled market-data refresh to Jenova_Profits sheets
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 275.25
$ws.Range("I2").Value = 275.25
$ws.Range("K2").Value = 275.25
$ws.Range("M2").Value = -162.25
$ws.Range("H11").Value = 416730.16
$ws.Range("I11").Value = 416730.16
$ws.Range("K11").Value = 416730.16
$ws.Range("M11").Value = -416590.16
$ws.Range("H64").Value = 9166.5
$ws.Range("H67").Value = 9166.5
$ws.Range("H74").Value = 5684.7646
$ws.Range("I74").Value = 3081.5
$ws.Range("K74").Value = 3081.5
$ws.Range("M74").Value = -2145.5
$ws.Range("H77").Value = 5684.7646
$ws.Range("I77").Value = 3081.5
$ws.Range("K77").Value = 15407.5
$ws.Range("M77").Value = -10727.5
$ws.Range("H112").Value = 3117.5
$ws.Range("J112").Value = 3171.3635
$ws.Range("L112").Value = 9514.0905
$ws.Range("N112").Value = -11730.0905
$ws.Range("H113").Value = 4997.143
$ws.Range("J113").Value = 4992.5
$ws.Range("L113").Value = 4992.5
$ws.Range("N113").Value = -11500.5
$ws.Range("H132").Value = 3023.5217
$ws.Range("I132").Value = 1798
$ws.Range("K132").Value = 5394
$ws.Range("M132").Value = -2864
$ws.Range("H137").Value = 4640.25
$ws.Range("I137").Value = 2964.4783
$ws.Range("J137").Value = 6475.619
$ws.Range("K137").Value = 8893.4349
$ws.Range("L137").Value = 19426.857
$ws.Range("M137").Value = -6343.4349
$ws.Range("N137").Value = -24526.857
$ws.Range("H138").Value = 5856.963
$ws.Range("I138").Value = 3632.5715
$ws.Range("K138").Value = 10897.7145
$ws.Range("M138").Value = -5757.7145

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 334632.66
$ws.Range("I2").Value = 334632.66
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 334632.66
$ws.Range("L2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("N2").Value = -334519.66
$ws.Range("H32").Value = 3784.05
$ws.Range("I32").Value = 3079.7114
$ws.Range("K32").Value = 3079.7114
$ws.Range("M32").Value = -2792.7114
$ws.Range("H39").Value = 36190.8
$ws.Range("I39").Value = 25850
$ws.Range("J39").Value = 77554
$ws.Range("K39").Value = 25850
$ws.Range("L39").Value = 77554
$ws.Range("M39").Value = -25330
$ws.Range("N39").Value = -78594
$ws.Range("H74").Value = 1508.2106
$ws.Range("I74").Value = 1303.931
$ws.Range("J74").Value = 2166.4443
$ws.Range("K74").Value = 1303.931
$ws.Range("L74").Value = 2166.4443
$ws.Range("M74").Value = -429.931
$ws.Range("N74").Value = -3914.4443
$ws.Range("H77").Value = 1508.2106
$ws.Range("I77").Value = 1303.931
$ws.Range("J77").Value = 2166.4443
$ws.Range("K77").Value = 6519.655000000001
$ws.Range("L77").Value = 10832.2215
$ws.Range("M77").Value = -2151.655000000001
$ws.Range("N77").Value = -19568.2215
$ws.Range("H112").Value = 56523.332
$ws.Range("J112").Value = 56523.332
$ws.Range("L112").Value = 56523.332
$ws.Range("N112").Value = -59477.332
$ws.Range("H116").Value = 334632.66
$ws.Range("I116").Value = 334632.66
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 334632.66
$ws.Range("L116").Value = 0
$ws.Range("M116").ClearContents()
$ws.Range("N116").Value = -332338.66
$ws.Range("H122").Value = 3163.6597
$ws.Range("I122").Value = 2084.1714
$ws.Range("K122").Value = 6252.514200000001
$ws.Range("M122").Value = -3802.514200000001
$ws.Range("H132").Value = 4888.4165
$ws.Range("I132").Value = 1475.5862
$ws.Range("J132").Value = 10097.474
$ws.Range("K132").Value = 4426.7586
$ws.Range("L132").Value = 30292.422
$ws.Range("M132").Value = -1896.7586
$ws.Range("N132").Value = -35352.422

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 334632.66
$ws.Range("I3").Value = 334632.66
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 334632.66
$ws.Range("L3").Value = 0
$ws.Range("M3").ClearContents()
$ws.Range("N3").Value = -334518.66
$ws.Range("H9").Value = 50000
$ws.Range("J9").Value = 50000
$ws.Range("L9").Value = 50000
$ws.Range("N9").Value = -50336
$ws.Range("H68").Value = 132647.5
$ws.Range("J68").Value = 132647.5
$ws.Range("L68").Value = 132647.5
$ws.Range("N68").Value = -134269.5
$ws.Range("H71").Value = 132647.5
$ws.Range("J71").Value = 132647.5
$ws.Range("L71").Value = 397942.5
$ws.Range("N71").Value = -406054.5
$ws.Range("H86").Value = 1423.9546
$ws.Range("I86").Value = 1214.3334
$ws.Range("K86").Value = 1214.3334
$ws.Range("M86").Value = -91.33339999999998
$ws.Range("H89").Value = 1423.9546
$ws.Range("I89").Value = 1214.3334
$ws.Range("K89").Value = 6071.666999999999
$ws.Range("M89").Value = -455.6669999999995
$ws.Range("H107").Value = 1527.3462
$ws.Range("I107").Value = 1443.9565
$ws.Range("K107").Value = 1443.9565
$ws.Range("M107").Value = 476.0435
$ws.Range("H134").Value = 3758.568
$ws.Range("I134").Value = 2467.5151
$ws.Range("K134").Value = 7402.5453
$ws.Range("M134").Value = -4867.5453

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H38").Value = 30000
$ws.Range("I38").Value = 30000
$ws.Range("K38").Value = 30000
$ws.Range("M38").Value = -29623
$ws.Range("H46").Value = 30000
$ws.Range("I46").Value = 30000
$ws.Range("K46").Value = 30000
$ws.Range("M46").Value = -29789

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 297.53845
$ws.Range("J12").Value = 342.63635
$ws.Range("L12").Value = 1027.90905
$ws.Range("N12").Value = -1373.90905
$ws.Range("H34").Value = 55271.047
$ws.Range("J34").Value = 96674.914
$ws.Range("L34").Value = 290024.742
$ws.Range("N34").Value = -290192.742
$ws.Range("H132").Value = 5087.5
$ws.Range("J132").Value = 5087.5
$ws.Range("L132").Value = 45787.5
$ws.Range("N132").Value = -50847.5
$ws.Range("H134").Value = 1778
$ws.Range("I134").Value = 1778
$ws.Range("K134").Value = 5334
$ws.Range("M134").Value = -264
$ws.Range("H137").Value = 3283.25
$ws.Range("I137").Value = 2377.6667
$ws.Range("J137").Value = 6000
$ws.Range("K137").Value = 7133.000100000001
$ws.Range("L137").Value = 18000
$ws.Range("M137").Value = -2033.000100000001
$ws.Range("N137").Value = -28200

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 135.875
$ws.Range("I2").Value = 164.92308
$ws.Range("K2").Value = 164.92308
$ws.Range("M2").Value = -51.92308
$ws.Range("H3").Value = 3130219.2
$ws.Range("J3").Value = 1673459.1
$ws.Range("L3").Value = 1673459.1
$ws.Range("N3").Value = -1673691.1
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("H70").Value = 40006490
$ws.Range("I70").Value = 5157.4
$ws.Range("J70").Value = 66674044
$ws.Range("K70").Value = 5157.4
$ws.Range("L70").Value = 66674044
$ws.Range("M70").Value = -4887.4
$ws.Range("N70").Value = -66674584
$ws.Range("H73").Value = 40006490
$ws.Range("I73").Value = 5157.4
$ws.Range("J73").Value = 66674044
$ws.Range("K73").Value = 5157.4
$ws.Range("L73").Value = 66674044
$ws.Range("M73").Value = -4221.4
$ws.Range("N73").Value = -66675916
$ws.Range("H95").Value = 19944
$ws.Range("J95").Value = 19944
$ws.Range("L95").Value = 19944
$ws.Range("N95").Value = -25436
$ws.Range("H97").Value = 9120.5
$ws.Range("I97").Value = 10874.7
$ws.Range("J97").Value = 349.5
$ws.Range("K97").Value = 10874.7
$ws.Range("L97").Value = 349.5
$ws.Range("M97").Value = -10378.7
$ws.Range("N97").Value = -1341.5
$ws.Range("H102").Value = 2415.4285
$ws.Range("I102").Value = 2415.4285
$ws.Range("K102").Value = 2415.4285
$ws.Range("M102").Value = -793.4285

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1258900.1
$ws.Range("I7").Value = 1437314.4
$ws.Range("K7").Value = 1437314.4
$ws.Range("M7").Value = -1437202.4
$ws.Range("H40").Value = 336111.03
$ws.Range("I40").Value = 347480.4
$ws.Range("K40").Value = 347480.4
$ws.Range("M40").Value = -347344.4
$ws.Range("H122").Value = 1175482.5
$ws.Range("I122").Value = 780564.7
$ws.Range("J122").Value = 2202268.5
$ws.Range("K122").Value = 2341694.1
$ws.Range("L122").Value = 6606805.5
$ws.Range("M122").Value = -2339244.1
$ws.Range("N122").Value = -6611705.5
$ws.Range("H126").Value = 1258900.1
$ws.Range("I126").Value = 1437314.4
$ws.Range("K126").Value = 4311943.199999999
$ws.Range("M126").Value = -4309473.199999999
$ws.Range("H132").Value = 4451.278
$ws.Range("I132").Value = 2921.1428
$ws.Range("J132").Value = 5425
$ws.Range("K132").Value = 8763.428400000001
$ws.Range("L132").Value = 16275
$ws.Range("M132").Value = -6233.428400000001
$ws.Range("N132").Value = -21335
$ws.Range("H136").Value = 3418.639
$ws.Range("I136").Value = 1947.2142
$ws.Range("J136").Value = 8568.625
$ws.Range("K136").Value = 5841.642599999999
$ws.Range("L136").Value = 25705.875
$ws.Range("M136").Value = -3291.642599999999
$ws.Range("N136").Value = -30805.875

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 49991.527
$ws.Range("J15").Value = 49991.527
$ws.Range("L15").Value = 49991.527
$ws.Range("N15").Value = -50567.527
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()
$ws.Range("H136").Value = 2478.6
$ws.Range("I136").Value = 2248.5417
$ws.Range("J136").Value = 8000
$ws.Range("K136").Value = 6745.625100000001
$ws.Range("L136").Value = 24000
$ws.Range("M136").Value = -4195.625100000001
$ws.Range("N136").Value = -29100

Write-Host "Applied scheduled Jenova_Profits update."
